# Natmi following Dr Hou advice
# Re-run with Ligand-/Receptor-expressing cell counts of 3 (was 1) for the
# Inha -> Acvr2a pairs; update the dependent average/total expression and
# specificity columns (G:H, M:T) for data rows 2-6 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => E, G, H, K, M, N, O, P, Q, R, S, T
$updates = @{
    2 = @{ E=3; G=0.5953293333333333; H=1.785988;          K=3; M=13.663472;          N=40.990416;          O=0.2163030132250088;  P=0.2163030132250088;  Q=8.134265676778664;  R=73.208391091008;   S=0.2163030132250088;  T=0.2163030132250088 }
    3 = @{ E=3; G=0.5953293333333333; H=1.785988;          K=3; M=30.94593233333333;  N=92.83779699999999;  O=0.4898973270305839;  P=0.4898973270305839;  Q=18.42302126538177;  R=165.807191388436;  S=0.4898973270305839;  T=0.4898973270305839 }
    4 = @{ E=3; G=0.5953293333333333; H=1.785988;          K=3; M=5.057967000000001;  N=15.173901;          O=0.08007141251452472; P=0.08007141251452471; Q=3.011156122132;     R=27.100405099188;   S=0.08007141251452472; T=0.08007141251452471 }
    5 = @{ E=3; G=0.5953293333333333; H=1.785988;          K=3; M=6.127540333333333;  N=18.382621;          O=0.09700356086343023; P=0.09700356086343022; Q=3.647904501616444;  R=32.831140514548;   S=0.09700356086343023; T=0.09700356086343022 }
    6 = @{ E=3; G=0.5953293333333333; H=1.785988;          K=3; M=7.373288333333332;  N=22.119865;          O=0.1167246863664523;  P=0.1167246863664523;  Q=4.389534827957776;  R=39.50581345161999; S=0.1167246863664523;  T=0.1167246863664523 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
    $ws.Range("H$row").Value = $vals.H
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("Q$row").Value = $vals.Q
    $ws.Range("R$row").Value = $vals.R
    $ws.Range("S$row").Value = $vals.S
    $ws.Range("T$row").Value = $vals.T
}
